{"js": "// Meeting 14 log line: \"09/22/2019 @ 2:00PM to 16:30PM\" -> \"09/22/2019 @ 2:00PM to 5:00PM\"\nconst meeting14 = context.document.body.search(\"09/22/2019 @ 2:00PM to 16:30PM\", { matchCase: true });\nmeeting14.load(\"text\");\nawait context.sync();\n\nif (meeting14.items.length > 0) {\n  meeting14.items[0].insertText(\"09/22/2019 @ 2:00PM to 5:00PM\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// Meeting 15 log line: \"09/22/2019 @ 17:00PM to 18:00PM\" -> \"09/22/2019 @ 5:00PM to 6:10PM\"\nconst meeting15 = context.document.body.search(\"09/22/2019 @ 17:00PM to 18:00PM\", { matchCase: true });\nmeeting15.load(\"text\");\nawait context.sync();\n\nif (meeting15.items.length > 0) {\n  meeting15.items[0].insertText(\"09/22/2019 @ 5:00PM to 6:10PM\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n# Meeting 14 log line: \"09/22/2019 @ 2:00PM to 16:30PM\" -> \"09/22/2019 @ 2:00PM to 5:00PM\"\n$find14 = $d.Content.Find\n$find14.ClearFormatting()\n$find14.Text = \"09/22/2019 @ 2:00PM to 16:30PM\"\n$find14.Replacement.ClearFormatting()\n$find14.Replacement.Text = \"09/22/2019 @ 2:00PM to 5:00PM\"\n$find14.Execute($find14.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find14.Replacement.Text, $wdReplaceAll) | Out-Null\n\n# Meeting 15 log line: \"09/22/2019 @ 17:00PM to 18:00PM\" -> \"09/22/2019 @ 5:00PM to 6:10PM\"\n$find15 = $d.Content.Find\n$find15.ClearFormatting()\n$find15.Text = \"09/22/2019 @ 17:00PM to 18:00PM\"\n$find15.Replacement.ClearFormatting()\n$find15.Replacement.Text = \"09/22/2019 @ 5:00PM to 6:10PM\"\n$find15.Execute($find15.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find15.Replacement.Text, $wdReplaceAll) | Out-Null\n"}
